$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal TEXT value (never let Excel
# auto-convert numeric-looking strings to actual numbers), then
# restore the cell to the default "Normal" style so no stray
# number-format/style is left behind.
function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "63.925.36"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "3.135.77"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "565.97"
$ws.Range("E5").Value = "  -0.58%  "
Set-TextValue $ws.Range("D6") "161.13"
$ws.Range("E6").Value = "  -4.50%  "
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue $ws.Range("D8") "0.568"
$ws.Range("E8").Value = "  -6.37%  "
Set-TextValue $ws.Range("D9") "0.115"
$ws.Range("E9").Value = "  -4.24%  "
Set-TextValue $ws.Range("D10") "6.55"
$ws.Range("E10").Value = "  -2.56%  "
Set-TextValue $ws.Range("D11") "0.378"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "3.692.65"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "64.095.14"
$ws.Range("E14").Value = "  -0.31%  "
Set-TextValue $ws.Range("D15") "24.82"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").Value = "3.145.60"
$ws.Range("E16").Value = "  -0.69%  "
Set-TextValue $ws.Range("D17") "0.0000153"
$ws.Range("E17").Value = "  -2.95%  "
Set-TextValue $ws.Range("D18") "398.51"
$ws.Range("E18").Value = "  -4.78%  "
Set-TextValue $ws.Range("D19") "12.51"
$ws.Range("E19").Value = "  -2.13%  "
Set-TextValue $ws.Range("D20") "5.18"
$ws.Range("E20").Value = "  -3.29%  "
Set-TextValue $ws.Range("D21") "7.06"
$ws.Range("E21").Value = "  +0.24%  "
Set-TextValue $ws.Range("D22") "5.85"
$ws.Range("E22").Value = "  +3.77%  "
$ws.Range("E23").Value = "  +0.05%  "
Set-TextValue $ws.Range("D24") "68.09"
$ws.Range("E24").Value = "  -3.02%  "
Set-TextValue $ws.Range("D25") "0.479"
$ws.Range("E25").Value = "  -2.19%  "
Set-TextValue $ws.Range("D26") "0.192"
$ws.Range("E26").Value = "  -4.85%  "
$ws.Range("D27").Value = "0.0₃0998"
$ws.Range("E27").Value = "  -5.46%  "
Set-TextValue $ws.Range("D28") "8.70"
$ws.Range("E28").Value = "  -0.58%  "
Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  +0.93%  "
Set-TextValue $ws.Range("D30") "1.79"
$ws.Range("E30").Value = "  -1.54%  "
Set-TextValue $ws.Range("D31") "20.95"
$ws.Range("E31").Value = "  -3.76%  "
Set-TextValue $ws.Range("D32") "6.19"
$ws.Range("E32").Value = "  -2.20%  "
Set-TextValue $ws.Range("D33") "4.76"
$ws.Range("E33").Value = "  -4.84%  "
Set-TextValue $ws.Range("D34") "156.32"
$ws.Range("E34").Value = "  +0.44%  "
Set-TextValue $ws.Range("D35") "1.10"
$ws.Range("E35").Value = "  -3.37%  "
Set-TextValue $ws.Range("D36") "1.31"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("D37").Value = "2.649.22"
$ws.Range("E37").Value = "  -1.92%  "
Set-TextValue $ws.Range("D38") "1.66"
$ws.Range("E38").Value = "  -1.96%  "
Set-TextValue $ws.Range("D39") "23.47"
$ws.Range("E39").Value = "  -4.15%  "
Set-TextValue $ws.Range("D40") "4.03"
$ws.Range("E40").Value = "  -3.39%  "
Set-TextValue $ws.Range("D41") "0.689"
$ws.Range("E41").Value = "  -2.54%  "
Set-TextValue $ws.Range("D42") "0.0608"
$ws.Range("E42").Value = "  -2.11%  "
Set-TextValue $ws.Range("D43") "5.41"
$ws.Range("E43").Value = "  -4.87%  "
$ws.Range("E44").Value = "  -2.65%  "
Set-TextValue $ws.Range("D45") "284.74"
$ws.Range("E45").Value = "  -3.46%  "
Set-TextValue $ws.Range("D46") "20.88"
$ws.Range("E46").Value = "  -4.04%  "
Set-TextValue $ws.Range("D47") "1.00"
$ws.Range("E47").Value = "  +0.01%  "
Set-TextValue $ws.Range("D48") "0.0972"
$ws.Range("E48").Value = "  -1.85%  "
Set-TextValue $ws.Range("D49") "10.49"
$ws.Range("E49").Value = "  +0.61%  "
Set-TextValue $ws.Range("D50") "1.86"
$ws.Range("E50").Value = "  -8.16%  "
Set-TextValue $ws.Range("D51") "5.62"
$ws.Range("E51").Value = "  -2.59%  "
